$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '67.856.54'
Set-TextValue 'E2' '  +0.75%  '
Set-TextValue 'D3' '3.485.33'
Set-TextValue 'E3' '  -0.09%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '594.40'
Set-TextValue 'E5' '  -0.43%  '
Set-TextValue 'D6' '182.41'
Set-TextValue 'E6' '  +3.16%  '
Set-TextValue 'E7' '  +4.95%  '
Set-TextValue 'E8' '  -0.06%  '
Set-TextValue 'D9' '3.481.82'
Set-TextValue 'E9' '  -0.16%  '
Set-TextValue 'D10' '0.141'
Set-TextValue 'E10' '  +7.56%  '
Set-TextValue 'D11' '7.01'
Set-TextValue 'E11' '  -1.27%  '
Set-TextValue 'D12' '0.430'
Set-TextValue 'E12' '  +0.96%  '
Set-TextValue 'D13' '4.081.96'
Set-TextValue 'E13' '  -0.15%  '
Set-TextValue 'D14' '32.12'
Set-TextValue 'E14' '  +1.51%  '
Set-TextValue 'E15' '  -1.06%  '
Set-TextValue 'D16' '67.813.74'
Set-TextValue 'E16' '  +0.71%  '
Set-TextValue 'D17' '0.0000178'
Set-TextValue 'E17' '  +0.77%  '
Set-TextValue 'D18' '3.485.29'
Set-TextValue 'E18' '  +0.13%  '
Set-TextValue 'D19' '6.21'
Set-TextValue 'E19' '  -0.66%  '
Set-TextValue 'D20' '14.15'
Set-TextValue 'E20' '  -1.88%  '
Set-TextValue 'D21' '395.46'
Set-TextValue 'E21' '  +1.67%  '
Set-TextValue 'D22' '7.97'
Set-TextValue 'E22' '  +0.57%  '
Set-TextValue 'D23' '5.83'
Set-TextValue 'E23' '  +1.95%  '
Set-TextValue 'B24' 'Polygon'
Set-TextValue 'C24' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D24' '0.540'
Set-TextValue 'E24' '  +0.82%  '
Set-TextValue 'B25' 'Dai'
Set-TextValue 'C25' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '0.999'
Set-TextValue 'E25' '  +0.00%  '
Set-TextValue 'D26' '72.23'
Set-TextValue 'E26' '  -0.95%  '
Set-TextValue 'D27' '0.0000123'
Set-TextValue 'E27' '  +0.65%  '
Set-TextValue 'D28' '10.37'
Set-TextValue 'E28' '  +0.91%  '
Set-TextValue 'E29' '  -0.81%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.33%  '
Set-TextValue 'D31' '6.14'
Set-TextValue 'E31' '  -0.26%  '
Set-TextValue 'D32' '1.41'
Set-TextValue 'E32' '  -1.20%  '
Set-TextValue 'E33' '  +0.78%  '
Set-TextValue 'D34' '23.66'
Set-TextValue 'E34' '  +0.25%  '
Set-TextValue 'D35' '7.35'
Set-TextValue 'E35' '  +0.76%  '
Set-TextValue 'E36' '  -0.02%  '
Set-TextValue 'E37' '  -3.42%  '
Set-TextValue 'D38' '161.60'
Set-TextValue 'E38' '  -1.27%  '
Set-TextValue 'D39' '0.896'
Set-TextValue 'E39' '  +2.84%  '
Set-TextValue 'D40' '2.89'
Set-TextValue 'E40' '  +11.94%  '
Set-TextValue 'E41' '  -2.99%  '
Set-TextValue 'D42' '4.70'
Set-TextValue 'E42' '  +1.20%  '
Set-TextValue 'D43' '6.76'
Set-TextValue 'E43' '  -3.36%  '
Set-TextValue 'D44' '26.26'
Set-TextValue 'E44' '  -0.79%  '
Set-TextValue 'D45' '0.0719'
Set-TextValue 'E45' '  -0.52%  '
Set-TextValue 'B46' 'Maker'
Set-TextValue 'C46' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D46' '2.749.71'
Set-TextValue 'E46' '  -2.38%  '
Set-TextValue 'B47' 'InjectiveProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D47' '26.32'
Set-TextValue 'E47' '  -3.38%  '
Set-TextValue 'D48' '41.63'
Set-TextValue 'E48' '  -1.44%  '
Set-TextValue 'E49' '  +0.12%  '
Set-TextValue 'D50' '329.34'
Set-TextValue 'E50' '  -3.64%  '
Set-TextValue 'D51' '1.05'
Set-TextValue 'E51' '  -2.27%  '
